$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.441.36'
$ws.Range('E2').Value = '  +2.38%  '
$ws.Range('D3').Value = '3.925.17'
$ws.Range('E3').Value = '  +4.11%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '470.79'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +10.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.734'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('E10').Value = '  +10.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000338'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +10.28%  '
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.550.52'
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.43'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.05'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '3.962.61'
$ws.Range('E16').Value = '  +4.07%  '
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.91'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('E19').Value = '  +4.14%  '
$ws.Range('D20').Value = '67.712.57'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.75'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +7.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.67'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.34'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.06'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.52%  '
$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '38.83'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.64%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.55'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +8.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.74'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.10'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.58'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '726.40'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.71'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.131'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.36%  '
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.05'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.62%  '
$ws.Range('E35').Value = '  +4.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.77'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('E37').Value = '  +20.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('E39').Value = '  -5.72%  '
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.04'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.06%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.336'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.54'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +5.27%  '
$ws.Range('E47').Value = '  +6.27%  '
$ws.Range('E48').Value = '  +1.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.18'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.43'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.89'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.25%  '
